$d = $word.ActiveDocument

# Remove the extra "M2Doc version mismatch" warning line that was appended
# after the "Could not activate TCPClientConnector[...]" text:
#   - the 4 spaces run that introduced it
#   - the orange "<---" marker run
#   - the "M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0" run
# while leaving everything else (including the following "    " run that
# precedes "demonstration") untouched.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = ""
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
